$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E, shifting the old "Warehouse" column to F
$ws.Columns("E:E").Insert()

# Set the new header
$ws.Range("E1").Value = "std_demand"

# Fill in the std_demand values for each data row
$ws.Range("E2").Value = 855
$ws.Range("E3").Value = 855
$ws.Range("E4").Value = 1045
$ws.Range("E5").Value = 1092.5
$ws.Range("E6").Value = 1377.5
$ws.Range("E7").Value = 1425
$ws.Range("E8").Value = 1615
$ws.Range("E9").Value = 1282.5
$ws.Range("E10").Value = 1045
$ws.Range("E11").Value = 712.5
$ws.Range("E12").Value = 807.5
$ws.Range("E13").Value = 950
$ws.Range("E14").Value = 997.5
$ws.Range("E15").Value = 1187.5
$ws.Range("E16").Value = 997.5
$ws.Range("E17").Value = 997.5
$ws.Range("E18").Value = 807.5
$ws.Range("E19").Value = 807.5
$ws.Range("E20").Value = 902.5
$ws.Range("E21").Value = 997.5
$ws.Range("E22").Value = 1187.5
$ws.Range("E23").Value = 1520
$ws.Range("E24").Value = 1425
$ws.Range("E25").Value = 1425
$ws.Range("E26").Value = 997.5
$ws.Range("E27").Value = 950
$ws.Range("E28").Value = 855
$ws.Range("E29").Value = 950
$ws.Range("E30").Value = 1092.5
$ws.Range("E31").Value = 950
$ws.Range("E32").Value = 902.5
$ws.Range("E33").Value = 807.5
$ws.Range("E34").Value = 997.5
$ws.Range("E35").Value = 1187.5
$ws.Range("E36").Value = 1187.5
$ws.Range("E37").Value = 1377.5
$ws.Range("E38").Value = 1140
$ws.Range("E39").Value = 1377.5
$ws.Range("E40").Value = 1472.5
$ws.Range("E41").Value = 1615
$ws.Range("E42").Value = 1710
$ws.Range("E43").Value = 1282.5
$ws.Range("E44").Value = 1330
$ws.Range("E45").Value = 997.5
$ws.Range("E46").Value = 1140
$ws.Range("E47").Value = 1140
$ws.Range("E48").Value = 1140
$ws.Range("E49").Value = 1282.5
$ws.Range("E50").Value = 1282.5
$ws.Range("E51").Value = 1472.5
$ws.Range("E52").Value = 1282.5
$ws.Range("E53").Value = 6127.5
$ws.Range("E54").Value = 6555
$ws.Range("E55").Value = 5985
$ws.Range("E56").Value = 4227.5
$ws.Range("E57").Value = 2683.75
$ws.Range("E58").Value = 1591.25
$ws.Range("E59").Value = 2018.75
$ws.Range("E60").Value = 1947.5
$ws.Range("E61").Value = 2232.5
$ws.Range("E62").Value = 2280
$ws.Range("E63").Value = 2470
$ws.Range("E64").Value = 3467.5
$ws.Range("E65").Value = 3325
$ws.Range("E66").Value = 3040
$ws.Range("E67").Value = 2232.5
$ws.Range("E68").Value = 2565
$ws.Range("E69").Value = 2802.5
$ws.Range("E70").Value = 3325
$ws.Range("E71").Value = 3230
$ws.Range("E72").Value = 2945
$ws.Range("E73").Value = 2185
$ws.Range("E74").Value = 2612.5
$ws.Range("E75").Value = 3562.5
$ws.Range("E76").Value = 4417.5
$ws.Range("E77").Value = 3847.5
$ws.Range("E78").Value = 3420
$ws.Range("E79").Value = 2565
$ws.Range("E80").Value = 2945
$ws.Range("E81").Value = 3895
$ws.Range("E82").Value = 4322.5
$ws.Range("E83").Value = 5130
$ws.Range("E84").Value = 5225
$ws.Range("E85").Value = 5842.5
$ws.Range("E86").Value = 5035
$ws.Range("E87").Value = 4180
$ws.Range("E88").Value = 5605
$ws.Range("E89").Value = 5225
$ws.Range("E90").Value = 4845
$ws.Range("E91").Value = 2232.5
$ws.Range("E92").Value = 2992.5
$ws.Range("E93").Value = 2850
$ws.Range("E94").Value = 3420
$ws.Range("E95").Value = 2137.5
$ws.Range("E96").Value = 3040
$ws.Range("E97").Value = 2802.5
$ws.Range("E98").Value = 3467.5
$ws.Range("E99").Value = 3040
$ws.Range("E100").Value = 2327.5
$ws.Range("E101").Value = 3372.5
$ws.Range("E102").Value = 4607.5
$ws.Range("E103").Value = 5320
$ws.Range("E104").Value = 3752.5
